$wb = $excel.ActiveWorkbook

# --- 3rd dose sheet: fill in zeros for the previously-blank weekly rows
#     (2021-07 through 2021-12, i.e. rows 2-32) and normalize a few
#     red-highlighted cells back to the regular number style. ---
$ws3 = $wb.Worksheets.Item("3rd dose")

$ws3.Range("B2:J32").Value = 0

# Rows 14 and 16, columns D:J were styled in red (s=3); reset them to the
# plain comma style (s=2) used by the rest of the sheet, same as copying
# format from a normal cell.
$ws3.Range("C14").Copy()
$ws3.Range("D14:J14").PasteSpecial(-4122)
$ws3.Range("C16").Copy()
$ws3.Range("D16:J16").PasteSpecial(-4122)

# D18 and D19 were individually styled red (s=3); reset just those cells.
$ws3.Range("E18").Copy()
$ws3.Range("D18").PasteSpecial(-4122)
$ws3.Range("E19").Copy()
$ws3.Range("D19").PasteSpecial(-4122)

$ws3.Range("B2:J32").Value = 0
$excel.CutCopyMode = 0

# Remove the old threaded comments sitting on D18/D19 (12-17 age group
# footnotes that no longer apply now that the rows are zeroed out).
$ws3.Range("D18").Comment.Delete()
$ws3.Range("D19").Comment.Delete()

# --- Update each sheet's scroll position / selection and make "3rd dose"
#     the active tab (previously "1st dose" was active). ---
$ws1 = $wb.Worksheets.Item("1st dose")
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 16
$ws1.Range("C3").Select()

$ws2 = $wb.Worksheets.Item("2nd dose")
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 15

$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 22
$ws3.Range("H27").Select()
